$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.260.76'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.666.88'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.78'
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5233'
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2677'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06345'
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.14'
$ws.Range("E10").Value = '  -2.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07733'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '1.678.83'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.444'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").Value = '1.891.92'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5500'
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Value = '0.0₅8248'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.14'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '26.298.29'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.676'
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '196.07'
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.18'
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.106'
$ws.Range("E23").Value = '  -3.95%  '
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.10'
$ws.Range("E25").Value = '  -2.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1248'
$ws.Range("E26").Value = '  -2.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.248'
$ws.Range("E27").Value = '  -2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.25'
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.414'
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05994'
$ws.Range("E30").Value = '  -3.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.287'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.616'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.313'
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.642'
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9857'
$ws.Range("E35").Value = '  -2.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.424'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.784'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5915'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.047'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01603'
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8608'
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  +0.34%  '
$ws.Range("D43").Value = '1.034.12'
$ws.Range("E43").Value = '  -3.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.37'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '1.805.47'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.67'
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈108'
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.017'
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.100'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.471'
$ws.Range("E51").Value = '  +1.30%  '
